{"js": "// Locate the paragraph that describes the \"SistemaCompatibleIndeterminado\"\n// case (it contains this distinctive phrase) and append the extra sentence\n// about the solution set condition to it, as a new run at its end.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"sistema de ecuaciones es compatible indeterminado\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"No se encontr\u00f3 el p\u00e1rrafo de descripci\u00f3n del SCI.\");\n}\n\ntarget.insertText(\n  \" La condici\u00f3n del conjunto soluci\u00f3n es S{(r-4;-2r+11;r)} con r un n\u00famero real.\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n", "ps1": "# Locate the paragraph describing the \"SistemaCompatibleIndeterminado\" case\n# (it contains this distinctive phrase) and append the extra sentence about\n# the solution set condition to the end of that paragraph, as a new run.\n$d = $word.ActiveDocument\n\n$marker = \"sistema de ecuaciones es compatible indeterminado\"\n$addition = \" La condici\u00f3n del conjunto soluci\u00f3n es S{(r-4;-2r+11;r)} con r un n\u00famero real.\"\n\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"*$marker*\") {\n        $target = $para\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"No se encontr\u00f3 el p\u00e1rrafo de descripci\u00f3n del SCI.\"\n}\n\n$target.Range.InsertAfter($addition)\n"}
